$d = $word.ActiveDocument

$replacements = @(
    @("630÷5=", "948÷3="),
    @("988÷4=", "779÷3="),
    @("299÷9=", "293÷4="),
    @("114÷4=", "328÷8="),
    @("677÷6=", "856÷6="),
    @("268÷6=", "825÷8="),
    @("573÷8=", "617÷9="),
    @("931÷7=", "600÷4="),
    @("943÷7=", "799÷4="),
    @("399÷9=", "856÷6="),
    @("108÷5=", "252÷6="),
    @("131÷4=", "657÷2="),
    @("120÷6=", "104÷7="),
    @("195÷8=", "154÷5="),
    @("890÷9=", "701÷7="),
    @("693÷2=", "474÷8="),
    @("630÷4=", "145÷7="),
    @("113÷7=", "142÷9="),
    @("610÷3=", "897÷6="),
    @("696÷7=", "992÷7="),
    @("465÷3=", "249÷5="),
    @("456÷8=", "392÷8="),
    @("765÷3=", "752÷6="),
    @("827÷8=", "484÷3="),
    @("911÷7=", "942÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
